$d = $word.ActiveDocument

# 1. "...elegant with on 21 lines..." -> "...elegant with only 21 lines..."
$d.Content.Find.Execute("elegant with on 21 lines", $true, $false, $false, $false, $false,
                         $true, 1, $false, "elegant with only 21 lines", 2) | Out-Null

# 2. "Four benchmarks were carried" -> "Three benchmarks were carried"
$d.Content.Find.Execute("Four benchmarks were carried", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Three benchmarks were carried", 2) | Out-Null

# 3. Remove the 64GB bulk key hashing clause
$d.Content.Find.Execute("the hashmap short key hashing, the 16MB bulk key hashing and the 64GB bulk key hashing.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false, "the hashmap short key hashing and the 16MB bulk key hashing.", 2) | Out-Null

# 4. Simplify repeat clause
$d.Content.Find.Execute(" and repeat 10000 times for short key hashing and 16MB hashing.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false, " and repeat 10000 times.", 2) | Out-Null

# 5. "XXH3 failed" -> "XXH3 fails"
$d.Content.Find.Execute("XXH3 failed two", $true, $false, $false, $false, $false,
                         $true, 1, $false, "XXH3 fails two", 2) | Out-Null

# 6. "Otma" + bookmark + "r Ertl" -> "Otmar Ertl" (merge text; bookmark relocated below)
$d.Content.Find.Execute("Otma", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Otmar", 2) | Out-Null
$d.Content.Find.Execute("r Ertl", $true, $false, $false, $false, $false,
                         $true, 1, $false, " Ertl", 2) | Out-Null

# 7. Move _GoBack bookmark from the Acknowledgements section to right after "XXH3 fails"
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmRange = $d.Content
$bmRange.Find.Execute("XXH3 fails") | Out-Null
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 8. Remove proofErr (spell-check) wrappers around "dumblob" and "ivte-ms" by
#    re-typing the surrounding text (fresh insertion carries no stale markers)
#    and deleting the old, marker-wrapped copy.
function Clean-ProofErr($searchText) {
    $r = $word.ActiveDocument.Content
    $r.Find.Execute($searchText) | Out-Null
    $text = $r.Text
    $origStart = $r.Start
    $origEnd = $r.End
    $insPoint = $r.Duplicate
    $insPoint.Collapse(1)
    $insPoint.Text = $text
    $len = $text.Length
    $toDelete = $word.ActiveDocument.Range($origStart + $len, $origEnd + $len)
    $toDelete.Delete()
}

Clean-ProofErr("paulie-g, dumblob, Yann")
Clean-ProofErr("Yann Collet, ivte-ms.")

Write-Output "done"
